# PYTHON_11_Functions Part 2.pptx
# Slide 30 ("Practice Questions"), Content Placeholder 2, 2nd paragraph:
#   "Write a function to multiply two numbers and return the product."
# becomes a 3-run paragraph:
#   "Write a function to multiply two numbers and return the product. "
#   "Or"  (red, FF0000)
#   " Build a calculator that adds, substracts, multiplies & divides two numbers."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(30)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# The paragraph we need to touch is the second one in this placeholder.
$target = $tr.Paragraphs(2, 1)

# First collapse the paragraph to a throwaway placeholder string. This
# drops the engine's memory of the old run boundaries so the text we set
# next does not get re-split at the position the old run used to end at.
$target.Text = "x"

# Re-fetch the (now 1-character) 2nd paragraph and set it to the full,
# final sentence in one shot so it starts life as a single run.
$target = $tr.Paragraphs(2, 1)
$newParagraphText = "Write a function to multiply two numbers and return the product. Or Build a calculator that adds, substracts, multiplies & divides two numbers."
$target.Text = $newParagraphText

# Work out where "Or" sits inside the whole placeholder's text so we can
# grab it with Characters(start, length) and recolor just that word.
$paraStart = $target.Start
$orOffset = $newParagraphText.IndexOf("Or")
$orStart = $paraStart + $orOffset

$orRun = $tr.Characters($orStart, 2)
$orRun.Font.Color.RGB = 255
